$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "X"
$ws.Range("B1").Value = "Y"

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 4

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 5

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 8

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 10

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 12

$ws.Range("A7").Value = 9
$ws.Range("B7").Value = 14

$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").Copy()
$ws.Range("A1:B7").PasteSpecial(-4122)

[void]$ws.Range("L9").Select()
